$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the two added columns (I = hasDiabetes, J = hasCancer)
$ws.Range("J1").Value = "hasCancer"
$ws.Range("I1").Value = "hasDiabetes"

$ws.Range("I2").Value = 12.9
$ws.Range("J2").Value = 8
$ws.Range("I3").Value = 10.1
$ws.Range("J3").Value = 6.9
$ws.Range("I4").Value = 12.8
$ws.Range("J4").Value = 6.8
$ws.Range("I5").Value = 10.3
$ws.Range("J5").Value = 6
$ws.Range("I6").Value = 7.3
$ws.Range("J6").Value = 6.1
$ws.Range("I7").Value = 9.1999999999999993
$ws.Range("J7").Value = 7.1
$ws.Range("I8").Value = 11.1
$ws.Range("J8").Value = 7.3
$ws.Range("I9").Value = 8.4
$ws.Range("J9").Value = 6.1
$ws.Range("I10").Value = 11.2
$ws.Range("J10").Value = 7.2
$ws.Range("I11").Value = 11.6
$ws.Range("J11").Value = 5.9
$ws.Range("I12").Value = 9.8000000000000007
$ws.Range("J12").Value = 5.8
$ws.Range("I13").Value = 7.6
$ws.Range("J13").Value = 6.8
$ws.Range("I14").Value = 10.1
$ws.Range("J14").Value = 5.3
$ws.Range("I15").Value = 10.7
$ws.Range("J15").Value = 6.3
$ws.Range("I16").Value = 9.5
$ws.Range("J16").Value = 6.8
$ws.Range("I17").Value = 10.3
$ws.Range("J17").Value = 6.8
$ws.Range("I18").Value = 12.5
$ws.Range("J18").Value = 8.1
$ws.Range("I19").Value = 11.3
$ws.Range("J19").Value = 6.5
$ws.Range("I20").Value = 10.1
$ws.Range("J20").Value = 5.8
$ws.Range("I21").Value = 9.6999999999999993
$ws.Range("J21").Value = 7.6
$ws.Range("I22").Value = 10.4
$ws.Range("J22").Value = 7.7
$ws.Range("I23").Value = 8.1
$ws.Range("J23").Value = 6.3
$ws.Range("I24").Value = 13
$ws.Range("J24").Value = 6.5
$ws.Range("I25").Value = 11.1
$ws.Range("J25").Value = 7.2
$ws.Range("I26").Value = 9.1999999999999993
$ws.Range("J26").Value = 6.1
$ws.Range("I27").Value = 9.6
$ws.Range("J27").Value = 5.9
$ws.Range("I28").Value = 9.1999999999999993
$ws.Range("J28").Value = 8
$ws.Range("I29").Value = 9.6999999999999993
$ws.Range("J29").Value = 6
$ws.Range("I30").Value = 11.5
$ws.Range("J30").Value = 6.9
$ws.Range("I31").Value = 10
$ws.Range("J31").Value = 5.8
$ws.Range("I32").Value = 10.8
$ws.Range("J32").Value = 6.2
$ws.Range("I33").Value = 11.7
$ws.Range("J33").Value = 6.4
$ws.Range("I34").Value = 12
$ws.Range("J34").Value = 6.2
$ws.Range("I35").Value = 9
$ws.Range("J35").Value = 7.9
$ws.Range("I36").Value = 11.2
$ws.Range("J36").Value = 6.9
$ws.Range("I37").Value = 9.5
$ws.Range("J37").Value = 6.8
$ws.Range("I38").Value = 12
$ws.Range("J38").Value = 7
$ws.Range("I39").Value = 13
$ws.Range("J39").Value = 7.4
$ws.Range("I40").Value = 11
$ws.Range("J40").Value = 5.3
$ws.Range("I41").Value = 7.1
$ws.Range("J41").Value = 5.5
$ws.Range("I42").Value = 9.6999999999999993
$ws.Range("J42").Value = 6.3
$ws.Range("I43").Value = 8.9
$ws.Range("J43").Value = 6.9
$ws.Range("I44").Value = 9
$ws.Range("J44").Value = 7
$ws.Range("I45").Value = 15.7
$ws.Range("J45").Value = 4.3

# Widen columns G (noHighSchoolGrad) and H (hasHealthCare) to fit the
# existing header text, matching the author's manual column-width tweak.
$ws.Range("G1").ColumnWidth = 10.72
$ws.Range("H1").ColumnWidth = 12.3

# Move the active selection to K1, mirroring the post-edit workbook state.
$ws.Range("K1").Select() | Out-Null
